$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K column (was previously Strike#) values for rows 2-10
$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
